# Auto-generated Excel COM-interop script to apply cryptos list update
# Commit message: "Updated cryptos list on Wed Jul 24 13:03:00 UTC 2024 with GitHub Actions"
#
# The sheet refreshes its "Price" (col D) and "Volume(1h)" (col E) columns
# with new scraped values. All of these cells are stored as TEXT in the
# workbook (e.g. "66.304.19", "  -0.19%  "), not numbers, so we briefly force
# NumberFormat="@" before assignment -- otherwise Excel would coerce strings
# like "2.00"/"1.00" into the numbers 2/1 and drop the literal formatting.
# The original cell style is captured/restored so no formatting drifts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $origStyle = $range.Style
    $range.NumberFormat = '@'
    $range.Value = $value
    $range.Style = $origStyle
}

Set-TextValue $ws.Range('D2') '66.251.26'
Set-TextValue $ws.Range('E2') '  -0.41%  '

Set-TextValue $ws.Range('D3') '3.448.96'
Set-TextValue $ws.Range('E3') '  -1.76%  '

Set-TextValue $ws.Range('E4') '  +0.03%  '

Set-TextValue $ws.Range('D5') '583.46'
Set-TextValue $ws.Range('E5') '  -0.34%  '

Set-TextValue $ws.Range('D6') '176.76'
Set-TextValue $ws.Range('E6') '  +0.63%  '

Set-TextValue $ws.Range('D7') '0.625'
Set-TextValue $ws.Range('E7') '  +4.97%  '

Set-TextValue $ws.Range('E8') '  +0.01%  '

Set-TextValue $ws.Range('D9') '3.447.16'
Set-TextValue $ws.Range('E9') '  -1.64%  '

Set-TextValue $ws.Range('D10') '0.132'
Set-TextValue $ws.Range('E10') '  -1.00%  '

Set-TextValue $ws.Range('E11') '  +1.32%  '

Set-TextValue $ws.Range('E12') '  -0.97%  '

Set-TextValue $ws.Range('D13') '4.045.83'
Set-TextValue $ws.Range('E13') '  -1.77%  '

Set-TextValue $ws.Range('E14') '  +1.30%  '

Set-TextValue $ws.Range('D15') '29.91'
Set-TextValue $ws.Range('E15') '  -1.72%  '

Set-TextValue $ws.Range('D16') '66.206.53'
Set-TextValue $ws.Range('E16') '  -0.40%  '

Set-TextValue $ws.Range('E17') '  -0.01%  '

Set-TextValue $ws.Range('D18') '3.440.35'
Set-TextValue $ws.Range('E18') '  -1.21%  '

Set-TextValue $ws.Range('E19') '  -0.77%  '

Set-TextValue $ws.Range('D20') '13.82'
Set-TextValue $ws.Range('E20') '  -0.65%  '

Set-TextValue $ws.Range('D21') '369.44'
Set-TextValue $ws.Range('E21') '  -2.66%  '

Set-TextValue $ws.Range('D22') '7.63'
Set-TextValue $ws.Range('E22') '  -3.27%  '

Set-TextValue $ws.Range('D23') '73.19'
Set-TextValue $ws.Range('E23') '  +1.47%  '

Set-TextValue $ws.Range('D24') '0.999'
Set-TextValue $ws.Range('E24') '  -0.16%  '

Set-TextValue $ws.Range('E25') '  -2.70%  '

Set-TextValue $ws.Range('E26') '  +3.72%  '

Set-TextValue $ws.Range('D27') '9.93'
Set-TextValue $ws.Range('E27') '  +0.40%  '

Set-TextValue $ws.Range('E28') '  +1.87%  '

Set-TextValue $ws.Range('E29') '  -0.03%  '

Set-TextValue $ws.Range('D30') '5.92'
Set-TextValue $ws.Range('E30') '  -0.08%  '

Set-TextValue $ws.Range('D31') '2.00'
Set-TextValue $ws.Range('E31') '  -0.79%  '

Set-TextValue $ws.Range('D32') '23.59'
Set-TextValue $ws.Range('E32') '  -3.79%  '

Set-TextValue $ws.Range('D33') '1.00'

Set-TextValue $ws.Range('D34') '7.05'
Set-TextValue $ws.Range('E34') '  -2.74%  '

Set-TextValue $ws.Range('E35') '  -5.04%  '

Set-TextValue $ws.Range('E36') '  -1.72%  '

Set-TextValue $ws.Range('D37') '161.60'
Set-TextValue $ws.Range('E37') '  +0.79%  '

Set-TextValue $ws.Range('D39') '27.85'
Set-TextValue $ws.Range('E39') '  -6.17%  '

Set-TextValue $ws.Range('E40') '  +0.61%  '

Set-TextValue $ws.Range('D41') '4.50'
Set-TextValue $ws.Range('E41') '  -0.52%  '

Set-TextValue $ws.Range('D42') '2.766.31'
Set-TextValue $ws.Range('E42') '  +2.67%  '

Set-TextValue $ws.Range('D43') '2.55'
Set-TextValue $ws.Range('E43') '  -0.36%  '

Set-TextValue $ws.Range('D44') '6.43'
Set-TextValue $ws.Range('E44') '  -0.36%  '

Set-TextValue $ws.Range('E45') '  -0.92%  '

Set-TextValue $ws.Range('D46') '25.12'
Set-TextValue $ws.Range('E46') '  +1.43%  '

Set-TextValue $ws.Range('D47') '339.75'
Set-TextValue $ws.Range('E47') '  +7.49%  '

Set-TextValue $ws.Range('D48') '39.95'
Set-TextValue $ws.Range('E48') '  -2.00%  '

Set-TextValue $ws.Range('D49') '0.0287'
Set-TextValue $ws.Range('E49') '  -1.72%  '

Set-TextValue $ws.Range('E50') '  +2.19%  '

Set-TextValue $ws.Range('D51') '0.991'
Set-TextValue $ws.Range('E51') '  -2.09%  '
